$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H112").Value = 1687.175
$ws.Range("J112").Value = 1764.5135
$ws.Range("L112").Value = 5293.5405
$ws.Range("N112").Value = -7509.5405
$ws.Range("H114").Value = 29851
$ws.Range("J114").Value = 29851
$ws.Range("L114").Value = 29851
$ws.Range("N114").Value = -38529
$ws.Range("H129").Value = 689.11
$ws.Range("I129").Value = 318.57144
$ws.Range("K129").Value = 955.71432
$ws.Range("M129").Value = 4044.28568
$ws.Range("H138").Value = 2847.889
$ws.Range("I138").Value = 1322.5834
$ws.Range("J138").Value = 3058.276
$ws.Range("K138").Value = 3967.7502
$ws.Range("L138").Value = 9174.828
$ws.Range("M138").Value = 1172.2498
$ws.Range("N138").Value = -19454.828

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 486.375
$ws.Range("I4").Value = 288.2
$ws.Range("J4").Value = 816.6667
$ws.Range("K4").Value = 288.2
$ws.Range("L4").Value = 816.6667
$ws.Range("M4").Value = -172.2
$ws.Range("N4").Value = -1048.6667
$ws.Range("H5").Value = 1000
$ws.Range("I5").Value = 1000
$ws.Range("K5").Value = 1000
$ws.Range("M5").Value = -888
$ws.Range("H6").Value = 759000
$ws.Range("I6").Value = 3000000
$ws.Range("J6").Value = 12000
$ws.Range("K6").Value = 3000000
$ws.Range("L6").Value = 12000
$ws.Range("M6").Value = -2999827
$ws.Range("N6").Value = -12346
$ws.Range("H23").Value = 11299
$ws.Range("I23").Value = 12497.5
$ws.Range("J23").Value = 10500
$ws.Range("K23").Value = 12497.5
$ws.Range("L23").Value = 10500
$ws.Range("M23").Value = -12238.5
$ws.Range("N23").Value = -11018
$ws.Range("H32").Value = 17566614
$ws.Range("I32").Value = 32268830
$ws.Range("J32").Value = 37050
$ws.Range("K32").Value = 32268830
$ws.Range("L32").Value = 37050
$ws.Range("M32").Value = -32268543
$ws.Range("N32").Value = -37624
$ws.Range("H37").Value = 17831.637
$ws.Range("I37").Value = 5034
$ws.Range("J37").Value = 19111.4
$ws.Range("K37").Value = 5034
$ws.Range("L37").Value = 19111.4
$ws.Range("M37").Value = -4761
$ws.Range("N37").Value = -19657.4
$ws.Range("H44").Value = 19571.428
$ws.Range("J44").Value = 19571.428
$ws.Range("L44").Value = 19571.428
$ws.Range("N44").Value = -20547.428
$ws.Range("H55").Value = 18015.143
$ws.Range("J55").Value = 18015.143
$ws.Range("L55").Value = 18015.143
$ws.Range("N55").Value = -18645.143
$ws.Range("H63").Value = 1740
$ws.Range("I63").Value = 1740
$ws.Range("K63").Value = 1740
$ws.Range("M63").Value = -1054
$ws.Range("H66").Value = 1740
$ws.Range("I66").Value = 1740
$ws.Range("K66").Value = 8700
$ws.Range("M66").Value = -5268
$ws.Range("H80").Value = 29642.857
$ws.Range("J80").Value = 29642.857
$ws.Range("L80").Value = 29642.857
$ws.Range("N80").Value = -31638.857
$ws.Range("H83").Value = 29642.857
$ws.Range("J83").Value = 29642.857
$ws.Range("L83").Value = 88928.571
$ws.Range("N83").Value = -98912.571

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 1000
$ws.Range("I4").Value = 1000
$ws.Range("K4").Value = 1000
$ws.Range("M4").Value = -885
$ws.Range("H15").Value = 34998.5
$ws.Range("J15").Value = 34998.5
$ws.Range("L15").Value = 34998.5
$ws.Range("N15").Value = -35452.5
$ws.Range("H35").Value = 17000
$ws.Range("J35").Value = 17000
$ws.Range("L35").Value = 17000
$ws.Range("N35").Value = -17620
$ws.Range("H82").Value = 16235.286
$ws.Range("I82").Value = 1011.75
$ws.Range("J82").Value = 36533.332
$ws.Range("K82").Value = 1011.75
$ws.Range("L82").Value = 36533.332
$ws.Range("M82").Value = -628.75
$ws.Range("N82").Value = -37299.332
$ws.Range("H85").Value = 16235.286
$ws.Range("I85").Value = 1011.75
$ws.Range("J85").Value = 36533.332
$ws.Range("K85").Value = 1011.75
$ws.Range("L85").Value = 36533.332
$ws.Range("M85").Value = 314.25
$ws.Range("N85").Value = -39185.332

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 207.33333
$ws.Range("I7").Value = 58.5
$ws.Range("J7").Value = 505
$ws.Range("K7").Value = 58.5
$ws.Range("L7").Value = 505
$ws.Range("M7").Value = 54.5
$ws.Range("N7").Value = -731

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 817.32074
$ws.Range("I131").Value = 431.1111
$ws.Range("J131").Value = 896.3182
$ws.Range("K131").Value = 1293.3333
$ws.Range("L131").Value = 2688.9546
$ws.Range("M131").Value = 3746.6667
$ws.Range("N131").Value = -12768.9546

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 135.85715
$ws.Range("I107").Value = 139.90909
$ws.Range("K107").Value = 139.90909
$ws.Range("M107").Value = 1780.09091
$ws.Range("H122").Value = 2730.1365
$ws.Range("I122").Value = 1259.7084
$ws.Range("J122").Value = 4494.65
$ws.Range("K122").Value = 3779.1252
$ws.Range("L122").Value = 13483.95
$ws.Range("M122").Value = -1329.1252
$ws.Range("N122").Value = -18383.95

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1092.7858
$ws.Range("I22").Value = 1249.75
$ws.Range("J22").Value = 1030
$ws.Range("K22").Value = 1249.75
$ws.Range("L22").Value = 1030
$ws.Range("M22").Value = -954.75
$ws.Range("N22").Value = -1620
$ws.Range("H27").Value = 1092.7858
$ws.Range("I27").Value = 1249.75
$ws.Range("J27").Value = 1030
$ws.Range("K27").Value = 1249.75
$ws.Range("L27").Value = 1030
$ws.Range("M27").Value = -1142.75
$ws.Range("N27").Value = -1244
$ws.Range("H68").Value = 2153.1765
$ws.Range("I68").Value = 1950.2858
$ws.Range("J68").Value = 3100
$ws.Range("K68").Value = 1950.2858
$ws.Range("L68").Value = 3100
$ws.Range("M68").Value = -1201.2858
$ws.Range("N68").Value = -4598
$ws.Range("H71").Value = 2153.1765
$ws.Range("I71").Value = 1950.2858
$ws.Range("J71").Value = 3100
$ws.Range("K71").Value = 9751.429
$ws.Range("L71").Value = 15500
$ws.Range("M71").Value = -6007.429
$ws.Range("N71").Value = -22988

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H127").Value = 47880
$ws.Range("J127").Value = 47880
$ws.Range("L127").Value = 47880
$ws.Range("N127").Value = -57800
